$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 8.6199999999999992
$ws.Range("E4").Value = 7.62
$ws.Range("E6").Value = 18
$ws.Range("E7").Value = 20.45
$ws.Range("E8").Value = 29.75
$ws.Range("E9").Value = 24.34
$ws.Range("E11").Value = 17.02
$ws.Range("E13").Value = 12.23
$ws.Range("E14").Value = 6.51
$ws.Range("E15").Value = 9.91
$ws.Range("E17").Value = 22.06
$ws.Range("E18").Value = 16.73
$ws.Range("E20").Value = 16.920000000000002
$ws.Range("E21").Value = 33.130000000000003
$ws.Range("E23").Value = 5.89
$ws.Range("E25").Value = 7.96
$ws.Range("E26").Value = 11.07
$ws.Range("E27").Value = 15.46
$ws.Range("E29").Value = 7.71
$ws.Range("E31").Value = 15.38
$ws.Range("E32").Value = 6.11
$ws.Range("E34").Value = 10.31
$ws.Range("E36").Value = 27.25
$ws.Range("E37").Value = 20.03
$ws.Range("E38").Value = 6.54
$ws.Range("E39").Value = 14.05
$ws.Range("E42").Value = 18.149999999999999
$ws.Range("E43").Value = 16.010000000000002
$ws.Range("E45").Value = 3.26
$ws.Range("E46").Value = 3.63
$ws.Range("E48").Value = 19.739999999999998
$ws.Range("E49").Value = 9.0399999999999991
$ws.Range("E50").Value = 4.43
$ws.Range("E52").Value = 6.4
$ws.Range("E53").Value = 12.96
$ws.Range("E55").Value = 10.36
$ws.Range("E56").Value = 6.27
$ws.Range("E58").Value = 20.82
$ws.Range("E60").Value = 8.43
$ws.Range("E61").Value = 9.83
$ws.Range("E62").Value = 4.4800000000000004
$ws.Range("E63").Value = 3.69
$ws.Range("E65").Value = 22.63
$ws.Range("E67").Value = 18.78
$ws.Range("E68").Value = 8.26
$ws.Range("E70").Value = 11.7
$ws.Range("E71").Value = 15.11
$ws.Range("E72").Value = 13.56

$ws.Range("E75").Select()

